$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output ("WARN: replace not found -> " + $find)
    }
}

# 1. Professional Summary: add AI safety research clause
Replace-Text "LLM inference efficiency and multi-agent systems." "LLM inference efficiency, multi-agent systems, and AI safety research (activation probing, sandbagging detection)."

# 2. 6 countries -> 7 countries (Data Analytics practice bullet)
Replace-Text "serving strategic enterprise clients across 6 countries." "serving strategic enterprise clients across 7 countries."

# 3. Standard Chartered summary bullet: drop "and 1200+ global users"
Replace-Text "serving 11 markets and 1200+ global users, delivering technical excellence" "serving 11 markets, delivering technical excellence"

# 4. Self-Service ML Platform: 6 months to 1 week -> months to weeks
Replace-Text "reduced model development time from 6 months to 1 week" "reduced model development time from months to weeks"

# 5. Architected enterprise-scale... drop trailing clause
Replace-Text "Architected enterprise-scale data solutions for Fortune 500 clients across APAC, designing scalable platforms with measurable business impact." "Architected enterprise-scale data solutions for Fortune 500 clients across APAC."

# 6. Various Companies -> Microsoft, Truckaurbus (Founder), UTU
Replace-Text "Various Companies" "Microsoft, Truckaurbus (Founder), UTU"

# 7. Role title line
Replace-Text "Software Engineering, Architecture and Technical Consulting Roles" "Software Engineering & Technical Leadership"

# 8. Progressively advanced -> Progressive advancement...
Replace-Text "Progressively advanced through roles in software development, systems integration, and technical consulting within financial services and algorithmic trading domains." "Progressive advancement through software engineering, entrepreneurship, and technical leadership across systems development, marketplace platforms, and payments infrastructure."

# 9. Modernized MarTech -> MarTech modernization - +30% customer acquisition
Replace-Text "Modernized MarTech infrastructure, driving 30% increase in customer acquisition" "MarTech modernization - +30% customer acquisition"

# 10. Engineered 5 high-performance data lakes -> Data lakes processing 1.2 PB/hour for Fortune 500 clients across APAC
Replace-Text "Engineered 5 high-performance data lakes processing 1.2 PB/hour, achieving 20% optimization" "Data lakes processing 1.2 PB/hour for Fortune 500 clients across APAC"

# 11. Built real-time fraud detection systems -> Real-time fraud detection systems - 60% reduction in false positives
Replace-Text "Built real-time fraud detection systems, reducing false positives by 60% and saving `$XM annually" "Real-time fraud detection systems - 60% reduction in false positives"

Write-Output "phase1-done"

function Delete-ParaContaining($needle) {
    $n = $d.Paragraphs.Count
    for ($i = $n; $i -ge 1; $i--) {
        $t = $d.Paragraphs($i).Range.Text
        if ($t -like "*$needle*") {
            $d.Paragraphs($i).Range.Delete()
            Write-Output ("deleted para containing: " + $needle)
            return
        }
    }
    Write-Output ("WARN: paragraph not found for deletion: " + $needle)
}

# Remove bullet: "Designed credit risk AI models..."
Delete-ParaContaining "Designed credit risk AI models integrating alternative data sources"

# Remove bullet: "Designed enterprise architectures supporting global Fortune 500 clients across APAC"
Delete-ParaContaining "Designed enterprise architectures supporting global Fortune 500 clients across APAC"

Write-Output "phase2-done"

# Remove the entire CatchMe project block (3 paragraphs):
#   "CatchMe - Intelligent Trust Engine (2025)"
#   "Google Technical Disclosures - Pending (APLS & Cascade Routing)"
#   "First of a kind, industry agnostic hybrid agentic AI decisioning system..."
Delete-ParaContaining "First of a kind, industry agnostic hybrid agentic AI decisioning system"
Delete-ParaContaining "Google Technical Disclosures - Pending (APLS & Cascade Routing)"
Delete-ParaContaining "CatchMe - Intelligent Trust Engine"

Write-Output "phase3-done"
